$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (shifts N:P -> O:Q),
# matching the author inserting a new "Variable Instalment" column
# on the repayment schedule sheet.
$ws.Columns("N:N").Insert()

# Match the original column's width (stored width "11") for the newly
# inserted column.
$ws.Columns("N:N").ColumnWidth = 10.166666666666666

# The user was working on the "Repayment schedule" tab (it becomes the
# active / selected sheet instead of "Transactions"), with the cursor
# left on R8 after the edit.
$ws.Activate()
$ws.Range("R8").Select()
